$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the sheet's gridline/row-col-header display (the COM round trip
# otherwise defaults these to hidden).
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.DisplayHeadings = $true

# 1. Update header row (row 1) from verbose labels to short field names.
$ws.Range("A1").Value = "email"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "password"
$ws.Range("D1").Value = "roles"
$ws.Range("E1").Value = "openai_api_key"

# 2. Replace the instructions comment on A1 with the new condensed text.
$newComment = @"
USER IMPORT TEMPLATE
Required fields:
• email - Must be unique and valid
• name - Full name (2-255 chars)
• password - Minimum 6 characters
Optional fields:
• roles - Comma-separated (e.g., instructor,admin)
• openai_api_key - User's OpenAI API key
Notes:
• Row 2 is an example - you can keep it or delete it
• Empty rows (no email) are ignored
• All imported users are verified
• Max file size: 5MB
"@
$ws.Range("A1").Comment.Text($newComment)

# 3. Remove the old "Instructions" block (rows 4-15) entirely, including the
#    merged cells that spanned A:E for each of those rows.
$ws.Rows("4:15").Delete()

# 4. Match the new selection left behind by the edit.
$ws.Range("E2").Select()
